$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Summoner Name / Summoner Role text on row 6
$ws.Range("G6").Value = "MyDogaN"
$ws.Range("H6").Value = "DUO_SUPPORT"

# Row 3 (A3 = 1)
$ws.Range("B3").Value = 2.732152492161565
$ws.Range("C3").Value = 4754.8
$ws.Range("D3").Value = 0.01058489950476475
$ws.Range("E3").Value = 18.2
$ws.Range("F3").Value = 286.2
$ws.Range("I3").Value = 0.1699525553653752
$ws.Range("J3").Value = 21.2
$ws.Range("K3").Value = 0.01215288164343084

# Row 4 (A4 = 2)
$ws.Range("B4").Value = 3.151435406698565
$ws.Range("C4").Value = 5269.2
$ws.Range("D4").Value = 0.02009569377990431
$ws.Range("E4").Value = 33.6
$ws.Range("F4").Value = 116
$ws.Range("I4").Value = 0.06937799043062201
$ws.Range("J4").Value = 2.6
$ws.Range("K4").Value = 0.001555023923444976

# Row 5 (A5 = 3)
$ws.Range("B5").Value = 6.383743300662566
$ws.Range("C5").Value = 12258.8
$ws.Range("D5").Value = 0.04023672888935409
$ws.Range("E5").Value = 81
$ws.Range("F5").Value = 253.6
$ws.Range("I5").Value = 0.13711238076999
$ws.Range("J5").Value = 14.8
$ws.Range("K5").Value = 0.007707535390238893

# Row 6 (A6 = 4)
$ws.Range("B6").Value = 1.901710291787398
$ws.Range("C6").Value = 4185.8
$ws.Range("D6").Value = 0.01349240180076153
$ws.Range("E6").Value = 29.6
$ws.Range("F6").Value = 204.4
$ws.Range("I6").Value = 0.09128055152368671
$ws.Range("J6").Value = 18.6
$ws.Range("K6").Value = 0.008167217339014521
